$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, shifting existing rows 31+ down by one.
$ws.Rows("31").Insert()

# Populate the newly inserted row 31 with the new product entry.
$ws.Range("A31").Value = "Salviette per la pulizia di scarpe da ginnastica usa e getta personalizzate in bambù non tessuto organico di alta qualità"
$ws.Range("B31").Value = "0,1733-0,3898 €"
$ws.Range("C31").Value = "Ordine minimo: 1.000 sacchi"
$ws.Range("D31").Value = "Zhejiang Furuisen Spunlaced Non-Wovens Co., Ltd."
$ws.Range("E31").Value = ""

# The last row (old row 49, now shifted to row 50) drops off the bottom of the
# page and is removed entirely.
$ws.Rows("50").Delete()
